$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.027.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.830.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  -0.76%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  +2.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07346"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8776"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07876"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.764.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.345"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.551"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.84%  "

$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008847"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.042.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.106"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.027.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("E26").Value = "  -0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.50"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.050"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.131"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08890"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.957"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7288"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.439"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.465"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.078"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01951"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05235"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.951"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.111"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5172"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1626"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.188"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.004"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.632"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06204"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.60%  "
